$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with rich-text runs) ---
$ws.Range("A8").Value = "Volume 31   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/19/2024  Through  8/25/2024"

# --- Column H width (cosmetic best-fit width) ---
$ws.Columns.Item(8).ColumnWidth = 5.4541602857143

# --- Text-forcing helper cells: set a reusable "@" text format once, then stamp it ---
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range("G14").NumberFormat = $ws.Range("F14").NumberFormat
$ws.Range("G14").Value = "0"
$ws.Range("H14").NumberFormat = $ws.Range("F14").NumberFormat
$ws.Range("H14").Value = "***.*"
$ws.Range("D31").NumberFormat = $ws.Range("F14").NumberFormat
$ws.Range("D31").Value = "0"
$ws.Range("E31").NumberFormat = $ws.Range("F14").NumberFormat
$ws.Range("E31").Value = "***.*"

# --- Cells switching from text placeholder to a real number: copy NumberFormat from a sibling cell that already has the right numeric style, then assign the value ---
$ws.Range("D29").NumberFormat = $ws.Range("F29").NumberFormat
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = $ws.Range("H29").NumberFormat
$ws.Range("E29").Value = -100
$ws.Range("D30").NumberFormat = $ws.Range("F30").NumberFormat
$ws.Range("D30").Value = 1
$ws.Range("E30").NumberFormat = $ws.Range("H30").NumberFormat
$ws.Range("E30").Value = -100
$ws.Range("F31").NumberFormat = $ws.Range("G31").NumberFormat
$ws.Range("F31").Value = 2
$ws.Range("C33").NumberFormat = $ws.Range("D33").NumberFormat
$ws.Range("C33").Value = 1

# --- Plain numeric value updates ---
# Row 14
$ws.Range("L14").Value = -47.058823529411
# Row 15
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = -25
$ws.Range("G15").Value = 11
$ws.Range("H15").Value = 72.727272727272
$ws.Range("I15").Value = 142
$ws.Range("J15").Value = 127
$ws.Range("K15").Value = 11.811023622047
$ws.Range("L15").Value = 10.9375
$ws.Range("M15").Value = 71.084337349397
$ws.Range("N15").Value = 8.396946564885
# Row 16
$ws.Range("C16").Value = 38
$ws.Range("D16").Value = 36
$ws.Range("E16").Value = 5.555555555555
$ws.Range("F16").Value = 152
$ws.Range("G16").Value = 170
$ws.Range("H16").Value = -10.588235294117
$ws.Range("I16").Value = 1441
$ws.Range("J16").Value = 1276
$ws.Range("K16").Value = 12.931034482758
$ws.Range("L16").Value = 27.635075287865
$ws.Range("M16").Value = 19.883527454242
$ws.Range("N16").Value = -74.917319408181
# Row 17
$ws.Range("C17").Value = 59
$ws.Range("D17").Value = 54
$ws.Range("E17").Value = 9.259259259259
$ws.Range("F17").Value = 260
$ws.Range("G17").Value = 243
$ws.Range("H17").Value = 6.995884773662
$ws.Range("I17").Value = 2197
$ws.Range("J17").Value = 1871
$ws.Range("K17").Value = 17.423837520042
$ws.Range("L17").Value = 34.290953545232
$ws.Range("M17").Value = 115.815324165029
$ws.Range("N17").Value = 14.965986394557
# Row 18
$ws.Range("C18").Value = 37
$ws.Range("D18").Value = 32
$ws.Range("E18").Value = 15.625
$ws.Range("F18").Value = 186
$ws.Range("G18").Value = 148
$ws.Range("H18").Value = 25.675675675675
$ws.Range("I18").Value = 1270
$ws.Range("J18").Value = 1288
$ws.Range("K18").Value = -1.39751552795
$ws.Range("L18").Value = 4.098360655737
$ws.Range("M18").Value = -25.513196480938
$ws.Range("N18").Value = -86.836650082918
# Row 19
$ws.Range("C19").Value = 144
$ws.Range("D19").Value = 148
$ws.Range("E19").Value = -2.702702702702
$ws.Range("F19").Value = 511
$ws.Range("G19").Value = 564
$ws.Range("H19").Value = -9.397163120567
$ws.Range("I19").Value = 4336
$ws.Range("J19").Value = 4461
$ws.Range("K19").Value = -2.802062317865
$ws.Range("L19").Value = -3.985828166519
$ws.Range("M19").Value = 71.45116646896
$ws.Range("N19").Value = -18.542175464963
# Row 20
$ws.Range("C20").Value = 61
$ws.Range("D20").Value = 59
$ws.Range("E20").Value = 3.389830508474
$ws.Range("F20").Value = 244
$ws.Range("G20").Value = 241
$ws.Range("H20").Value = 1.244813278008
$ws.Range("I20").Value = 1695
$ws.Range("J20").Value = 1634
$ws.Range("K20").Value = 3.733170134638
$ws.Range("L20").Value = 52.702702702702
$ws.Range("M20").Value = 46.373056994818
$ws.Range("N20").Value = -88.811142649679
# Row 21
$ws.Range("C21").Value = 342
$ws.Range("D21").Value = 333
$ws.Range("E21").Value = 2.702702702702
$ws.Range("F21").Value = 1372
$ws.Range("G21").Value = 1377
$ws.Range("H21").Value = -0.363108206245
$ws.Range("I21").Value = 11090
$ws.Range("J21").Value = 10672
$ws.Range("K21").Value = 3.916791604197
$ws.Range("L21").Value = 13.673636736367
$ws.Range("M21").Value = 43.764583873476
$ws.Range("N21").Value = -70.810412444397
# Row 22
$ws.Range("C22").Value = 10
$ws.Range("E22").Value = 150
$ws.Range("F22").Value = 29
$ws.Range("H22").Value = 26.086956521739
$ws.Range("I22").Value = 204
$ws.Range("J22").Value = 237
$ws.Range("K22").Value = -13.924050632911
$ws.Range("L22").Value = 13.333333333333
$ws.Range("M22").Value = 67.213114754098
# Row 23
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 17
$ws.Range("G23").Value = 27
$ws.Range("H23").Value = -37.037037037037
$ws.Range("I23").Value = 170
$ws.Range("J23").Value = 179
$ws.Range("K23").Value = -5.027932960893
$ws.Range("L23").Value = 3.658536585365
$ws.Range("M23").Value = 53.153153153153
# Row 24
$ws.Range("C24").Value = 328
$ws.Range("D24").Value = 348
$ws.Range("E24").Value = -5.747126436781
$ws.Range("F24").Value = 1289
$ws.Range("G24").Value = 1292
$ws.Range("H24").Value = -0.232198142414
$ws.Range("I24").Value = 10564
$ws.Range("J24").Value = 10046
$ws.Range("K24").Value = 5.156281106908
$ws.Range("L24").Value = 7.270511779041
$ws.Range("M24").Value = 74.958595561444
# Row 25
$ws.Range("C25").Value = 212
$ws.Range("D25").Value = 182
$ws.Range("E25").Value = 16.483516483516
$ws.Range("F25").Value = 787
$ws.Range("G25").Value = 715
$ws.Range("H25").Value = 10.06993006993
$ws.Range("I25").Value = 6513
$ws.Range("J25").Value = 5604
$ws.Range("K25").Value = 16.220556745182
$ws.Range("L25").Value = 28.486881041625
# Row 26
$ws.Range("C26").Value = 104
$ws.Range("D26").Value = 104
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 539
$ws.Range("G26").Value = 408
$ws.Range("H26").Value = 32.107843137254
$ws.Range("I26").Value = 4154
$ws.Range("J26").Value = 3394
$ws.Range("K26").Value = 22.392457277548
$ws.Range("L26").Value = 27.540681608842
$ws.Range("M26").Value = 32.461734693877
# Row 27
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 7
$ws.Range("E27").Value = -42.857142857142
$ws.Range("F27").Value = 26
$ws.Range("G27").Value = 20
$ws.Range("H27").Value = 30
$ws.Range("I27").Value = 211
$ws.Range("J27").Value = 201
$ws.Range("K27").Value = 4.975124378109
$ws.Range("L27").Value = 11.052631578947
# Row 28
$ws.Range("C28").Value = 11
$ws.Range("D28").Value = 17
$ws.Range("E28").Value = -35.294117647058
$ws.Range("G28").Value = 52
$ws.Range("H28").Value = -19.230769230769
$ws.Range("I28").Value = 400
$ws.Range("J28").Value = 473
$ws.Range("K28").Value = -15.433403805496
$ws.Range("L28").Value = -1.477832512315
# Row 29
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -25
$ws.Range("J29").Value = 46
$ws.Range("K29").Value = -65.217391304347
$ws.Range("L29").Value = -64.444444444444
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -90.857142857142
# Row 30
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -50
$ws.Range("J30").Value = 41
$ws.Range("K30").Value = -68.292682926829
$ws.Range("L30").Value = -68.292682926829
$ws.Range("M30").Value = -53.571428571428
$ws.Range("N30").Value = -91.772151898734
# Row 31
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 41
$ws.Range("K31").Value = -21.153846153846
$ws.Range("L31").Value = -10.869565217391
# Row 33
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = -50
$ws.Range("F33").Value = 3
$ws.Range("G33").Value = 8
$ws.Range("H33").Value = -62.5
$ws.Range("I33").Value = 32
$ws.Range("J33").Value = 33
$ws.Range("K33").Value = -3.030303030303
$ws.Range("L33").Value = 39.130434782608
